$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row values (row 1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 values
$ws.Range("B2").Value = 12.416459387778305
$ws.Range("C2").Value = 11.644600955605554
$ws.Range("D2").Value = 12.868442019214903
$ws.Range("E2").Value = 12.648640548210338

# Update row 3 values
$ws.Range("B3").Value = 12.59734379604013
$ws.Range("C3").Value = 10.669031043082921
$ws.Range("D3").Value = 14.317480407274433
$ws.Range("E3").Value = 11.513577455305288

# Update the selection to match the new selected range
$ws.Range("B1:E3").Select()
